$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# User re-opened the workbook (on a new server/machine), selected cell B3
# and corrected its value from 2600 to 2500.
$cell = $ws.Range("B3")
$cell.Select()
$cell.Value = 2500
